$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Jalen Brunson"
$ws.Range("B2").Value = "PG"
$ws.Range("C2").Value = "New York Knicks"

$ws.Range("A3").Value = "Devin Booker"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Phoenix Suns"

$ws.Range("A4").Value = "Norman Powell"
$ws.Range("B4").Value = "SG,SF"
$ws.Range("C4").Value = "LA Clippers"

$ws.Range("A5").Value = "D'Angelo Russell"
$ws.Range("B5").Value = "PG"
$ws.Range("C5").Value = "Brooklyn Nets"

$ws.Range("A6").Value = "Kawhi Leonard"
$ws.Range("B6").Value = "SG,SF,PF"
$ws.Range("C6").Value = "LA Clippers"

$ws.Range("A7").Value = "LeBron James"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Los Angeles Lakers"

$ws.Range("A8").Value = "Jimmy Butler"
$ws.Range("B8").Value = "SF,PF"
$ws.Range("C8").Value = "Miami Heat"

$ws.Range("A9").Value = "Desmond Bane"
$ws.Range("B9").Value = "SG,SF"
$ws.Range("C9").Value = "Memphis Grizzlies"

$ws.Range("A10").Value = "Walker Kessler"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "Utah Jazz"

$ws.Range("A11").Value = "Alexandre Sarr"
$ws.Range("B11").Value = "PF,C"
$ws.Range("C11").Value = "Washington Wizards"

$ws.Range("A12").Value = "Bilal Coulibaly"
$ws.Range("B12").Value = "SG,SF"
$ws.Range("C12").Value = "Washington Wizards"

$ws.Range("A13").Value = "Vit Krejci"
$ws.Range("B13").Value = "SG,SF"
$ws.Range("C13").Value = "Atlanta Hawks"

$ws.Range("A14").Value = "Jalen Williams"
$ws.Range("B14").Value = "SG,SF,PF,C"
$ws.Range("C14").Value = "Oklahoma City Thunder"

$ws.Range("A15").Value = "Scoot Henderson"
$ws.Range("B15").Value = "PG"
$ws.Range("C15").Value = "Portland Trail Blazers"

$ws.Range("A16").Value = "Myles Turner"
$ws.Range("B16").Value = "C"
$ws.Range("C16").Value = "Indiana Pacers"

$ws.Range("A17").Value = "Immanuel Quickley"
$ws.Range("B17").Value = "PG,SG"
$ws.Range("C17").Value = "Toronto Raptors"

$ws.Range("A18").Value = "Brandon Ingram"
$ws.Range("B18").Value = "SG,SF,PF"
$ws.Range("C18").Value = "New Orleans Pelicans"

$ws.Range("A19").Value = "Trae Young"
$ws.Range("B19").Value = "PG"
$ws.Range("C19").Value = "Atlanta Hawks"
